$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.418.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.824.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5107'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3919'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07644'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.105'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.262'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.502'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.822.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001094'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06671'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.137'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.444.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.263'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.033.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.382'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.106'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1086'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.627'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.656'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07033'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2203'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02314'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.812'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.149'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6230'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.168'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.388'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5876'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.706'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.83'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06911'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.23%  '
